$wb = $excel.ActiveWorkbook

# --- Insert new visible sheet "test-error-catch" right before the hidden
#     "__footings__" sheet (mirrors footings writing a new test sheet,
#     then re-appending the footings report sheet after it). ---
$footings = $wb.Worksheets.Item("__footings__")
$errCatch = $wb.Worksheets.Add($footings)
$errCatch.Name = "test-error-catch"

$errCatch.Cells.Item(2,2).Value = "key"
$errCatch.Cells.Item(2,3).Value = "{'k1': '1'}"
$errCatch.Cells.Item(3,2).Value = "error_type"
$errCatch.Cells.Item(3,3).Value = "TypeError"
$errCatch.Cells.Item(4,2).Value = "error_value"
$errCatch.Cells.Item(4,3).Value = "(`"__init__() missing 1 required keyword-only argument: 'k2'`",)"
$errCatch.Cells.Item(5,2).Value = "error_stacktrace"
$errCatch.Cells.Item(5,3).Value = "['  File `"/home/dustintindall/anaconda3/lib/python3.7/site-packages/footings/parallel_tools/base.py`", line 107, in wrapper`n    ret = model(**model_kwargs).run()`n']"

# --- Append the footings-report rows describing the new test-error-catch
#     sheet to the bottom of the (now shifted) "__footings__" sheet.
#     Re-fetch by name: inserting a sheet shifts index-based handles. ---
$footings = $wb.Worksheets.Item("__footings__")
$footings.Cells.Item(52,1).Value = "test-error-catch"
$footings.Cells.Item(52,3).Value = "/key/"
$footings.Cells.Item(52,4).Value = "KEY"
$footings.Cells.Item(52,6).Value = "<class 'str'>"
$footings.Cells.Item(52,8).Value = 2
$footings.Cells.Item(52,9).Value = 2
$footings.Cells.Item(52,10).Value = 2
$footings.Cells.Item(52,11).Value = 2

$footings.Cells.Item(53,1).Value = "test-error-catch"
$footings.Cells.Item(53,3).Value = "/key/"
$footings.Cells.Item(53,4).Value = "VALUE"
$footings.Cells.Item(53,6).Value = "<class 'str'>"
$footings.Cells.Item(53,8).Value = 2
$footings.Cells.Item(53,9).Value = 3
$footings.Cells.Item(53,10).Value = 2
$footings.Cells.Item(53,11).Value = 3

$footings.Cells.Item(54,1).Value = "test-error-catch"
$footings.Cells.Item(54,3).Value = "/error_type/"
$footings.Cells.Item(54,4).Value = "KEY"
$footings.Cells.Item(54,6).Value = "<class 'str'>"
$footings.Cells.Item(54,8).Value = 3
$footings.Cells.Item(54,9).Value = 2
$footings.Cells.Item(54,10).Value = 3
$footings.Cells.Item(54,11).Value = 2

$footings.Cells.Item(55,1).Value = "test-error-catch"
$footings.Cells.Item(55,3).Value = "/error_type/"
$footings.Cells.Item(55,4).Value = "VALUE"
$footings.Cells.Item(55,6).Value = "<class 'str'>"
$footings.Cells.Item(55,8).Value = 3
$footings.Cells.Item(55,9).Value = 3
$footings.Cells.Item(55,10).Value = 3
$footings.Cells.Item(55,11).Value = 3

$footings.Cells.Item(56,1).Value = "test-error-catch"
$footings.Cells.Item(56,3).Value = "/error_value/"
$footings.Cells.Item(56,4).Value = "KEY"
$footings.Cells.Item(56,6).Value = "<class 'str'>"
$footings.Cells.Item(56,8).Value = 4
$footings.Cells.Item(56,9).Value = 2
$footings.Cells.Item(56,10).Value = 4
$footings.Cells.Item(56,11).Value = 2

$footings.Cells.Item(57,1).Value = "test-error-catch"
$footings.Cells.Item(57,3).Value = "/error_value/"
$footings.Cells.Item(57,4).Value = "VALUE"
$footings.Cells.Item(57,6).Value = "<class 'str'>"
$footings.Cells.Item(57,8).Value = 4
$footings.Cells.Item(57,9).Value = 3
$footings.Cells.Item(57,10).Value = 4
$footings.Cells.Item(57,11).Value = 3

$footings.Cells.Item(58,1).Value = "test-error-catch"
$footings.Cells.Item(58,3).Value = "/error_stacktrace/"
$footings.Cells.Item(58,4).Value = "KEY"
$footings.Cells.Item(58,6).Value = "<class 'str'>"
$footings.Cells.Item(58,8).Value = 5
$footings.Cells.Item(58,9).Value = 2
$footings.Cells.Item(58,10).Value = 5
$footings.Cells.Item(58,11).Value = 2

$footings.Cells.Item(59,1).Value = "test-error-catch"
$footings.Cells.Item(59,3).Value = "/error_stacktrace/"
$footings.Cells.Item(59,4).Value = "VALUE"
$footings.Cells.Item(59,6).Value = "<class 'str'>"
$footings.Cells.Item(59,8).Value = 5
$footings.Cells.Item(59,9).Value = 3
$footings.Cells.Item(59,10).Value = 5
$footings.Cells.Item(59,11).Value = 3
